# Applies cryptos.xlsx price/volume/coin updates described in the commit
# "Updated cryptos list on Sat Apr 22 11:00:01 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that are plain text (coin names, links, percent strings with
# padding, and multi-dot price strings) can be written directly: Excel will
# not misinterpret them as numbers.
$textValues = [ordered]@{
    'D2' = '27.282.68'
    'E2' = '  -2.74%  '
    'D3' = '1.852.25'
    'E3' = '  -3.37%  '
    'E4' = '  -0.08%  '
    'E5' = '  -1.58%  '
    'E6' = '  -0.04%  '
    'E7' = '  -3.21%  '
    'E8' = '  -3.65%  '
    'E9' = '  -9.24%  '
    'E10' = '  -5.99%  '
    'E11' = '  -3.06%  '
    'E12' = '  -3.61%  '
    'D13' = '1.869.30'
    'E13' = '  -2.82%  '
    'E14' = '  -2.75%  '
    'E15' = '  -4.27%  '
    'E16' = '  +0.14%  '
    'B17' = 'Litecoin'
    'C17' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'E17' = '  -4.63%  '
    'B18' = 'TRON'
    'C18' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'E18' = '  +0.26%  '
    'E19' = '  -3.97%  '
    'E20' = '  -4.55%  '
    'E21' = '  +0.04%  '
    'E22' = '  -4.07%  '
    'D23' = '27.285.68'
    'E23' = '  -2.81%  '
    'E24' = '  -4.40%  '
    'E25' = '  +0.13%  '
    'D26' = '2.086.45'
    'E26' = '  -2.96%  '
    'E27' = '  +0.06%  '
    'E28' = '  -0.55%  '
    'E29' = '  -3.85%  '
    'E30' = '  -4.37%  '
    'E31' = '  -1.78%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'E32' = '  -2.65%  '
    'B33' = 'ImmutableX'
    'C33' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'E33' = '  -3.01%  '
    'E34' = '  -0.49%  '
    'E35' = '  -1.45%  '
    'E36' = '  -5.09%  '
    'E38' = '  -3.62%  '
    'E39' = '  -1.09%  '
    'E40' = '  -9.64%  '
    'E41' = '  -0.07%  '
    'E42' = '  -3.50%  '
    'E43' = '  -1.05%  '
    'E44' = '  -7.93%  '
    'E45' = '  -2.20%  '
    'E46' = '  -4.55%  '
    'E47' = '  -6.22%  '
    'E48' = '  -2.86%  '
    'E49' = '  -5.78%  '
    'E50' = '  -1.47%  '
    'E51' = '  -1.84%  '
}
foreach ($cell in $textValues.Keys) {
    $ws.Range($cell).Value = $textValues[$cell]
}

# Cell values that look like plain decimal numbers ("1.002", "325.25", ...)
# need to be forced to text so Excel stores them the same way the source
# workbook does (as strings, not numbers) - otherwise trailing zeros and
# the General number format would corrupt them. We temporarily switch the
# cell to a text number format, assign the literal string, then restore the
# cell style so no stray formatting is left behind.
$numericLookingTextValues = [ordered]@{
    'D4' = '1.002'
    'D5' = '325.25'
    'D7' = '0.4556'
    'D8' = '0.3885'
    'D9' = '48.27'
    'D10' = '0.07906'
    'D11' = '1.013'
    'D12' = '21.34'
    'D14' = '5.905'
    'D15' = '7.145'
    'D16' = '1.004'
    'D17' = '85.83'
    'D18' = '0.06597'
    'D19' = '0.00001024'
    'D20' = '17.20'
    'D22' = '5.489'
    'D24' = '10.84'
    'D25' = '2.291'
    'D27' = '154.08'
    'D28' = '19.90'
    'D29' = '2.057'
    'D30' = '5.453'
    'D31' = '121.20'
    'D32' = '0.09335'
    'D33' = '0.9429'
    'D34' = '1.439'
    'D35' = '3.588'
    'D36' = '5.253'
    'D37' = '0.06026'
    'D40' = '8.064'
    'D41' = '1.001'
    'D42' = '0.5918'
    'D43' = '0.1883'
    'D44' = '10.15'
    'D45' = '1.278'
    'D46' = '0.5592'
    'D47' = '12.01'
    'D48' = '3.380'
    'D49' = '1.909'
    'D51' = '107.80'
}
foreach ($cell in $numericLookingTextValues.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $numericLookingTextValues[$cell]
    $range.Style = "Normal"
}
